# Insert 9 new rows of historical data (2019-11-18 .. 2019-11-28) right
# after the existing row 574 (2019-11-15), shifting all following rows
# down by 9 (the old row 575 "2019-11-29" becomes row 584, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 blank rows starting at row 575 (pushes existing rows down).
$ws.Range("A575:A583").EntireRow.Insert()

# New rows data: row, timestamp, date, open, high, low, close, vol
$newRows = @(
    @(575, 1574035200, "2019-11-18", 2.44, 2.46, 2.44, 2.46, 154500),
    @(576, 1574121600, "2019-11-19", 2.46, 2.46, 2.44, 2.46, 226400),
    @(577, 1574208000, "2019-11-20", 2.46, 2.46, 2.44, 2.45, 280700),
    @(578, 1574294400, "2019-11-21", 2.45, 2.45, 2.42, 2.44, 1962300),
    @(579, 1574380800, "2019-11-22", 2.43, 2.44, 2.43, 2.44, 2505700),
    @(580, 1574640000, "2019-11-25", 2.43, 2.44, 2.41, 2.41, 399200),
    @(581, 1574726400, "2019-11-26", 2.41, 2.44, 2.41, 2.42, 271300),
    @(582, 1574812800, "2019-11-27", 2.42, 2.43, 2.41, 2.42, 66800),
    @(583, 1574899200, "2019-11-28", 2.42, 2.42, 2.3,  2.34, 1227700)
)

foreach ($row in $newRows) {
    $r = $row[0]

    # Force the id/name/date columns to be stored as text, matching the
    # rest of the sheet (they would otherwise be auto-converted to
    # numbers/dates by Excel's smart-typing).
    $ws.Range("B$r`:D$r").NumberFormat = "@"

    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = "5284"
    $ws.Range("D$r").Value = "LCTITAN"
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
    $ws.Range("H$r").Value = $row[6]
    $ws.Range("I$r").Value = $row[7]

    # Restore default (Normal) styling so the new cells look the same
    # as the rest of the data rows (no explicit style index).
    $ws.Range("B$r`:D$r").Style = "Normal"
}
